# Complete functional test protocol
#
# "Rectangle 23" (inside "Group 37") is dragged to the right by 378894 EMU.
# PowerPoint glues "Elbow Connector 24" to it (stCxn id="24"), so the
# connector's left edge follows while its right edge (glued to "Rectangle 1")
# stays put, shrinking its width by the same amount. The parent group's
# transform (off/ext/chOff/chExt on X) is recomputed by PowerPoint from the
# non-connector member shapes only, which also shifts its left edge right and
# shrinks its width by 378894 while the right edge stays fixed.
#
# This COM surface does not expose a "recompute group bounds" primitive or a
# way to set chOff/chExt directly, and moving a shape inside an existing
# group does not touch the parent's stored off/ext/chOff/chExt. The only
# operation observed to recompute off/ext/chOff/chExt together (keeping the
# 1:1 child scale) is ShapeRange.Group(), which derives the new box from the
# *current* Left/Top/Width/Height of every shape being grouped (including
# connectors, unlike real PowerPoint's manual-drag recompute).
#
# So: ungroup, move "Rectangle 23" to its final spot, temporarily collapse
# every connector shape onto a point inside the target box (so it can't
# inflate the computed bounds), regroup (which now yields exactly the target
# off/ext/chOff/chExt), then restore every connector to its real final
# geometry (5 unchanged, 1 resized) -- which leaves the freshly computed
# group box untouched, since child moves don't perturb it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# Index every (formerly grouped) shape by its PowerPoint shape Id.
$range = $grp.Ungroup()
$byId = @{}
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $byId[$shp.Id] = $shp
}

$connectorIds = @(4, 12, 16, 17, 25, 37)

# Each connector's true current (point-precise, exact-EMU-roundtrip) geometry.
$origLeft   = @{ 4 = 215.94166564941406; 12 = 215.2101593017578;  16 = 284.25347900390625; 17 = 228.67283630371094;  25 = 113.6415786743164;    37 = 335.8317565917969 }
$origTop    = @{ 4 = 138.55039978027344; 12 = 205.4117431640625;  16 = 189.30496215820312; 17 = 253.55615234375;     25 = 186.35284423828125;   37 = 292.46356201171875 }
$origWidth  = @{ 4 = 0.7314173579216003; 12 = 0.7314173579216003; 16 = 115.79544067382812; 17 = 59.482208251953125;  25 = 55.447166442871094;   37 = 0.7314173579216003 }
$origHeight = @{ 4 = 28.743623733520508; 12 = 28.743623733520508; 16 = 14.286378860473633; 17 = 19.848583221435547;  25 = 0.00007874015864217654; 37 = 28.743623733520508 }

# Move "Rectangle 23" (id 24) to its final position (x += 378894 EMU; y is
# unchanged).
$rect23 = $byId[24]
$rect23.Left = 79.85662078857422

# Collapse every connector onto a point that sits well inside the target
# group box so none of them can stretch the bounds Group() is about to
# compute.
foreach ($cid in $connectorIds) {
    $c = $byId[$cid]
    $c.Left = 79.85662078857422
    $c.Top = 174.1580352783203
    $c.Width = 0.00007874015748031496
    $c.Height = 0.00007874015748031496
}

# Regroup everything; with the connectors collapsed this reproduces exactly
# the target off/ext/chOff/chExt for the group (x: 1014179 / 3907729 EMU).
$newgrp = $range.Group()
$newgrp.Name = "Group 37"

# Restore the five untouched connectors to their exact original geometry.
foreach ($cid in @(4, 12, 16, 17, 37)) {
    $c = $byId[$cid]
    $c.Left = $origLeft[$cid]
    $c.Top = $origTop[$cid]
    $c.Width = $origWidth[$cid]
    $c.Height = $origHeight[$cid]
}

# "Elbow Connector 24" (id 25) follows "Rectangle 23": left edge moves to
# 1822142 EMU, width shrinks to 325285 EMU; y/height stay at their original
# values.
$conn24 = $byId[25]
$conn24.Top = $origTop[25]
$conn24.Height = $origHeight[25]
$conn24.Left = 143.4757537841797
$conn24.Width = 25.612993240356445
